# Applies the cryptos-list refresh described by the commit diff: per-row
# Price (D) / Volume(1h) (E) text updates, plus four coin swaps where two
# adjacent rows traded ranks (so B/C/D/E all change on those rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All these sheet cells are stored as plain text (inline strings) in the
# workbook, including ones that look like plain numbers (e.g. "7.95" or
# "4.30"). Force text format before writing so Excel does not silently
# coerce them to numeric values (which would also drop trailing zeros,
# e.g. turning "4.30" into 4.3).
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D9",
    "D11",
    "D14",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D36",
    "D39",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D47",
    "D48",
    "D51",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "60.939.28"
$ws.Range("E2").Value = "  -2.62%  "
$ws.Range("D3").Value = "3.355.63"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "566.69"
$ws.Range("E5").Value = "  -1.94%  "
$ws.Range("D6").Value = "147.15"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Value = "7.95"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").Value = "0.415"
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("D12").Value = "3.940.20"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "27.96"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "3.369.63"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "60.971.21"
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").Value = "6.31"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").Value = "14.42"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").Value = "8.92"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").Value = "376.35"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").Value = "0.561"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "75.15"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "3.505.78"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("E26").Value = "  -6.36%  "
$ws.Range("D27").Value = "0.175"
$ws.Range("E27").Value = "  -3.42%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "7.39"
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.08"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "7.69"
$ws.Range("E32").Value = "  -3.56%  "
$ws.Range("D33").Value = "22.85"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("D35").Value = "5.34"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").Value = "169.39"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  -3.63%  "
$ws.Range("E38").Value = "  -2.40%  "
$ws.Range("D39").Value = "28.87"
$ws.Range("E39").Value = "  -9.48%  "
$ws.Range("D40").Value = "3.392.08"
$ws.Range("E40").Value = "  -2.30%  "
$ws.Range("D41").Value = "0.0751"
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("D42").Value = "0.759"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").Value = "4.30"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.61"
$ws.Range("E44").Value = "  -5.16%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").Value = "1.14"
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("D46").Value = "2.492.53"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").Value = "22.65"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").Value = "6.67"
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").Value = "0.813"
$ws.Range("E51").Value = "  +0.07%  "
